# Update Metrics: inserted metrics for lr_cv_std_mix_features
#
# 1) The deck's "today" date placeholder (an auto-updating a:fld
#    type="datetime1") is cached in every slide layout + the slide
#    master + notes master + handout master. The original author's
#    save rolled the cached date from 20/06/2022 to 21/06/2022 - update
#    every one of those cached placeholders to match.
# 2) Slide 12 has the EXPERIMENTAL RESULTS table; the Logistic
#    Regression / "Mix - with standardization" cell was still blank -
#    fill it in with the 0.974 metric, matching the styling already
#    used by its sibling cells (Avenir Next LT Pro).

$p = $ppt.ActivePresentation

$newDate = "21/06/2022"

# --- Slide master date placeholder -----------------------------------
$sm = $p.SlideMaster
$sm.Shapes.Item(3).TextFrame.TextRange.Text = $newDate

# --- Notes master date placeholder ------------------------------------
$nm = $p.NotesMaster
$nm.Shapes.Item(2).TextFrame.TextRange.Text = $newDate

# --- Handout master date placeholder ----------------------------------
$hm = $p.HandoutMaster
$hm.Shapes.Item(2).TextFrame.TextRange.Text = $newDate

# --- Every slide layout's date placeholder -----------------------------
# CustomLayouts.Item(N) <-> ppt/slideLayouts/slideLayoutN.xml ; the date
# placeholder shape index differs per layout.
$layoutDateShapeIndex = @{
    1  = 4
    2  = 4
    3  = 4
    4  = 5
    5  = 7
    6  = 1
    7  = 1
    8  = 5
    9  = 4
    10 = 4
    11 = 4
}

for ($li = 1; $li -le $sm.CustomLayouts.Count; $li++) {
    $layout = $sm.CustomLayouts.Item($li)
    $shapeIdx = $layoutDateShapeIndex[$li]
    $layout.Shapes.Item($shapeIdx).TextFrame.TextRange.Text = $newDate
}

# --- Slide 12: fill in the missing metric cell --------------------------
$slide12 = $p.Slides.Item(12)
$tableShape = $slide12.Shapes.Item(3)
$table = $tableShape.Table

# Row 3 = "Logistic Regression", Column 6 = "Mix / with standardization"
$cell = $table.Cell(3, 6)
$cellRange = $cell.Shape.TextFrame.TextRange
$cellRange.Text = "0.974"
$cellRange.Font.Name = "Avenir Next LT Pro"
